$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format risky numeric-looking text cells as Text so Excel keeps them as strings
$textCells = @("D4","D5","D7","D8","D9","D10","D11","D13","D14","D15","D17","D18","D19","D20","D22","D25","D26","D28","D29","D31","D33","D34","D35","D37","D38","D39","D40","D41","D42","D43","D44","D48","D49","D50","D51","D45","D46")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply new values
$ws.Range("D2").Value = "26.831.52"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.813.88"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "307.99"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.4328"
$ws.Range("E7").Value = "  +2.47%  "
$ws.Range("D8").Value = "0.3712"
$ws.Range("E8").Value = "  +3.07%  "
$ws.Range("D9").Value = "0.07258"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "0.8669"
$ws.Range("E10").Value = "  +2.74%  "
$ws.Range("D11").Value = "20.91"
$ws.Range("E11").Value = "  +3.11%  "
$ws.Range("D12").Value = "1.930.49"
$ws.Range("E12").Value = "  +2.24%  "
$ws.Range("D13").Value = "6.651"
$ws.Range("E13").Value = "  +4.44%  "
$ws.Range("D14").Value = "5.361"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").Value = "0.06923"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "80.61"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "0.000008937"
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "15.22"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").Value = "26.857.67"
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("D22").Value = "5.216"
$ws.Range("E22").Value = "  +2.58%  "
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").Value = "2.162.91"
$ws.Range("E24").Value = "  +3.73%  "
$ws.Range("D25").Value = "153.93"
$ws.Range("D26").Value = "1.871"
$ws.Range("E26").Value = "  -4.33%  "
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("D28").Value = "5.226"
$ws.Range("E28").Value = "  +3.97%  "
$ws.Range("D29").Value = "1.906"
$ws.Range("E29").Value = "  +15.23%  "
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("D31").Value = "0.08941"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("E32").Value = "  +3.70%  "
$ws.Range("D33").Value = "1.173"
$ws.Range("E33").Value = "  +7.27%  "
$ws.Range("D34").Value = "4.439"
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("D35").Value = "2.811"
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("D37").Value = "1.134"
$ws.Range("E37").Value = "  +5.08%  "
$ws.Range("D38").Value = "0.05235"
$ws.Range("E38").Value = "  +1.59%  "
$ws.Range("D39").Value = "0.01926"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("D40").Value = "0.5090"
$ws.Range("E40").Value = "  +2.22%  "
$ws.Range("D41").Value = "0.1651"
$ws.Range("E41").Value = "  +1.15%  "
$ws.Range("D42").Value = "2.681"
$ws.Range("E42").Value = "  +2.52%  "
$ws.Range("D43").Value = "6.561"
$ws.Range("E43").Value = "  +9.86%  "
$ws.Range("D44").Value = "8.312"
$ws.Range("E44").Value = "  +2.67%  "
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").Value = "1.655"
$ws.Range("E48").Value = "  +3.29%  "
$ws.Range("D49").Value = "0.4586"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").Value = "0.06279"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").Value = "1.823"
$ws.Range("E51").Value = "  +5.50%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "10.38"
$ws.Range("E45").Value = "  +1.64%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "106.76"
$ws.Range("E46").Value = "  +1.58%  "

# Restore default (no explicit) formatting on cells we forced to text,
# so the workbook's style table matches the original structure.
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
